$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Model Accuracy (-1.0, 1.0, 1.0)"
# Add new columns C:G (Market threshold, Market min, Market max, Recall,
# Precision) and update column B accuracy values.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Model Accuracy (-1.0, 1.0, 1.0)")

# Header row - copy the existing "Accuracy (%)" header formatting (bold,
# centered, thin border) from B1 onto the new header cells, then set text.
$ws1.Range("B1").Copy()
$ws1.Range("C1:G1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws1.Range("C1").Value = "Market threshold"
$ws1.Range("D1").Value = "Market min"
$ws1.Range("E1").Value = "Market max"
$ws1.Range("F1").Value = "Recall"
$ws1.Range("G1").Value = "Precision"

# Row 2 - TOTALENERGIES SE
$ws1.Range("B2").Value = 67.1760391198044
$ws1.Range("C2").Value = 0.05450546436368681
$ws1.Range("D2").Value = -15.55441
$ws1.Range("E2").Value = 15.06418
$ws1.Range("F2").Value = 0
$ws1.Range("G2").Value = 0

# Row 3 - FMC CORP
$ws1.Range("B3").Value = 39.85330073349633
$ws1.Range("C3").Value = 0.009583939973006913
$ws1.Range("D3").Value = -19.35264
$ws1.Range("E3").Value = 13.70093
$ws1.Range("F3").Value = 0
$ws1.Range("G3").Value = 0

# Row 4 - BP PLC
$ws1.Range("B4").Value = 92.72616136919315
$ws1.Range("C4").Value = 0.04158117063764853
$ws1.Range("D4").Value = -18.75314
$ws1.Range("E4").Value = 23.33066
$ws1.Range("F4").Value = 0
$ws1.Range("G4").Value = 0

# Row 5 - STORA ENSO
$ws1.Range("B5").Value = 83.00733496332519
$ws1.Range("C5").Value = 0.02983403801513819
$ws1.Range("D5").Value = -12.78028
$ws1.Range("E5").Value = 12.42348
$ws1.Range("F5").Value = 0
$ws1.Range("G5").Value = 0

# Row 6 - BHP GROUP
$ws1.Range("B6").Value = 96.14914425427872
$ws1.Range("C6").Value = 0.08368817696170747
$ws1.Range("D6").Value = -16.47904
$ws1.Range("E6").Value = 14.94325
$ws1.Range("F6").Value = 0
$ws1.Range("G6").Value = 0

# ---------------------------------------------------------------------------
# Sheet 2: "Confusion Matrix TOTALENERGIES SE (-1.0, 1.0, 1.0)"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Confusion Matrix TOTALENERGIES SE (-1.0, 1.0, 1.0)")
$ws2.Range("B3").Value = 9
$ws2.Range("C3").Value = 1099
$ws2.Range("D3").Value = 10

# ---------------------------------------------------------------------------
# Sheet 3: "Confusion Matrix FMC CORP (-1.0, 1.0, 1.0)"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Confusion Matrix FMC CORP (-1.0, 1.0, 1.0)")
$ws3.Range("B3").Value = 373
$ws3.Range("C3").Value = 652
$ws3.Range("D3").Value = 356

# ---------------------------------------------------------------------------
# Sheet 4: "Confusion Matrix BP PLC (-1.0, 1.0, 1.0)"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Confusion Matrix BP PLC (-1.0, 1.0, 1.0)")
$ws4.Range("B3").Value = 40
$ws4.Range("C3").Value = 1517
$ws4.Range("D3").Value = 42

# ---------------------------------------------------------------------------
# Sheet 5: "Confusion Matrix STORA ENSO (-1.0, 1.0, 1.0)"
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Confusion Matrix STORA ENSO (-1.0, 1.0, 1.0)")
$ws5.Range("B3").Value = 110
$ws5.Range("C3").Value = 1358
$ws5.Range("D3").Value = 107

# ---------------------------------------------------------------------------
# Sheet 6: "Confusion Matrix BHP GROUP (-1.0, 1.0, 1.0)"
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Confusion Matrix BHP GROUP (-1.0, 1.0, 1.0)")
$ws6.Range("B3").Value = 4
$ws6.Range("C3").Value = 1573
$ws6.Range("D3").Value = 3
